$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values would otherwise be
# auto-coerced to numbers by Excel (so they stay text, matching original inlineStr type).
$textCells = @("D5","D6","D9","D10","D11","D19","D21","D23","D24","D25","D26","D28","D29","D30","D31","D32","D34","D35","D37","D39","D40","D41","D42","D43","D44","D49","D51")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply cell value updates
$ws.Range('D2').Value = '50.073.04'
$ws.Range('E2').Value = '  +4.27%  '
$ws.Range('D3').Value = '2.659.75'
$ws.Range('E3').Value = '  +7.01%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '114.12'
$ws.Range('E5').Value = '  +8.23%  '
$ws.Range('D6').Value = '326.68'
$ws.Range('E6').Value = '  +2.95%  '
$ws.Range('E7').Value = '  +2.09%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '0.558'
$ws.Range('E9').Value = '  +3.93%  '
$ws.Range('D10').Value = '41.27'
$ws.Range('E10').Value = '  +6.21%  '
$ws.Range('D11').Value = '20.17'
$ws.Range('E11').Value = '  +0.10%  '
$ws.Range('E12').Value = '  +3.14%  '
$ws.Range('E13').Value = '  +0.49%  '
$ws.Range('E14').Value = '  +4.22%  '
$ws.Range('D15').Value = '3.077.55'
$ws.Range('D16').Value = '2.662.49'
$ws.Range('E16').Value = '  +6.94%  '
$ws.Range('E17').Value = '  +6.23%  '
$ws.Range('D18').Value = '49.988.50'
$ws.Range('E18').Value = '  +4.28%  '
$ws.Range('D19').Value = '13.30'
$ws.Range('E19').Value = '  +4.65%  '
$ws.Range('E20').Value = '  +3.71%  '
$ws.Range('D21').Value = '2.94'
$ws.Range('E21').Value = '  -1.44%  '
$ws.Range('D22').Value = '0.0₃0961'
$ws.Range('E22').Value = '  +3.42%  '
$ws.Range('D23').Value = '72.60'
$ws.Range('E23').Value = '  +2.34%  '
$ws.Range('D24').Value = '278.47'
$ws.Range('E24').Value = '  +2.12%  '
$ws.Range('D25').Value = '2.61'
$ws.Range('E25').Value = '  +3.97%  '
$ws.Range('D26').Value = '26.91'
$ws.Range('E26').Value = '  +4.85%  '
$ws.Range('E27').Value = '  +0.10%  '
$ws.Range('D28').Value = '10.03'
$ws.Range('E28').Value = '  +3.49%  '
$ws.Range('D29').Value = '2.24'
$ws.Range('E29').Value = '  -1.86%  '
$ws.Range('D30').Value = '36.75'
$ws.Range('E30').Value = '  +6.46%  '
$ws.Range('D31').Value = '0.143'
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('D32').Value = '50.26'
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('E33').Value = '  +4.74%  '
$ws.Range('D34').Value = '19.73'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('D35').Value = '0.0818'
$ws.Range('E35').Value = '  +6.21%  '
$ws.Range('E36').Value = '  -0.21%  '
$ws.Range('D37').Value = '5.05'
$ws.Range('E37').Value = '  +10.70%  '
$ws.Range('E38').Value = '  +7.78%  '
$ws.Range('D39').Value = '3.12'
$ws.Range('E39').Value = '  +9.15%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D40').Value = '125.13'
$ws.Range('E40').Value = '  +2.05%  '
$ws.Range('B41').Value = 'Stellar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D41').Value = '0.113'
$ws.Range('E41').Value = '  +2.79%  '
$ws.Range('D42').Value = '22.38'
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('D43').Value = '2.23'
$ws.Range('E43').Value = '  +0.40%  '
$ws.Range('D44').Value = '0.0320'
$ws.Range('E44').Value = '  +5.35%  '
$ws.Range('D45').Value = '2.110.41'
$ws.Range('E45').Value = '  +5.51%  '
$ws.Range('E46').Value = '  +5.67%  '
$ws.Range('E47').Value = '  +13.50%  '
$ws.Range('E48').Value = '  +5.42%  '
$ws.Range('D49').Value = '9.11'
$ws.Range('E49').Value = '  +2.34%  '
$ws.Range('E50').Value = '  +3.94%  '
$ws.Range('D51').Value = '59.77'
$ws.Range('E51').Value = '  +5.86%  '
